# Apply cell value updates for cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.039.53'
$ws.Range('E2').Value = '  +2.33%  '
$ws.Range('D3').Value = '1.655.92'
$ws.Range('E3').Value = '  +2.87%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'215.30"
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('E6').Value = '  +2.40%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +2.30%  '
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('E10').Value = '  +4.71%  '
$ws.Range('D12').Value = '1.890.28'
$ws.Range('E12').Value = '  +2.95%  '
$ws.Range('D13').Value = '1.643.93'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('E14').Value = '  +1.88%  '
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('D16').Value = "'65.26"
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('D17').Value = '27.045.26'
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('D18').Value = "'236.57"
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0732'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = "'7.74"
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = "'4.43"
$ws.Range('E22').Value = '  +3.49%  '
$ws.Range('D23').Value = "'9.35"
$ws.Range('E23').Value = '  +3.95%  '
$ws.Range('E24').Value = '  +2.95%  '
$ws.Range('D25').Value = "'145.64"
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('D26').Value = "'7.11"
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').Value = "'15.85"
$ws.Range('E28').Value = '  +2.17%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = "'0.0497"
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('D32').Value = '1.558.75'
$ws.Range('E32').Value = '  +3.70%  '
$ws.Range('E33').Value = '  +2.99%  '
$ws.Range('D34').Value = "'3.08"
$ws.Range('E34').Value = '  +4.54%  '
$ws.Range('D35').Value = "'1.61"
$ws.Range('E35').Value = '  +8.38%  '
$ws.Range('D36').Value = "'2.42"
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').Value = "'0.584"
$ws.Range('E37').Value = '  +3.03%  '
$ws.Range('D38').Value = "'0.903"
$ws.Range('E38').Value = '  +9.33%  '
$ws.Range('E39').Value = '  +2.54%  '
$ws.Range('D40').Value = "'6.00"
$ws.Range('E40').Value = '  +2.92%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Value = "'2.26"
$ws.Range('E42').Value = '  +3.19%  '
$ws.Range('D43').Value = "'65.66"
$ws.Range('E43').Value = '  +7.47%  '
$ws.Range('D44').Value = '1.797.07'
$ws.Range('E44').Value = '  +2.78%  '
$ws.Range('D45').Value = "'0.776"
$ws.Range('E45').Value = '  +1.93%  '
$ws.Range('D46').Value = "'0.914"
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('D47').Value = "'90.40"
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').Value = '0.0₆0104'
$ws.Range('E48').Value = '  +11.52%  '
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('E50').Value = '  +2.21%  '
$ws.Range('E51').Value = '  +0.82%  '

Write-Output "Applied 82 cell updates"
